$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL): plain string assignment ---
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"

# --- Numeric-looking text columns (Price / Volume%): force text storage ---
$cellsToFormat = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $cellsToFormat) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "306.44"
$ws.Range("E2").Value = "-3.55%"
$ws.Range("D3").Value = "37.12"
$ws.Range("E3").Value = "-6.66%"
$ws.Range("D4").Value = "5.065"
$ws.Range("E4").Value = "-1.48%"
$ws.Range("D5").Value = "0.07722"
$ws.Range("E5").Value = "-6.07%"
$ws.Range("D6").Value = "4.340"
$ws.Range("E6").Value = "0.64%"
$ws.Range("D7").Value = "1.896"
$ws.Range("E7").Value = "-7.46%"
$ws.Range("D8").Value = "8.189"
$ws.Range("E8").Value = "-2.32%"
$ws.Range("E9").Value = "-5.72%"
$ws.Range("D10").Value = "0.9189"
$ws.Range("E10").Value = "-2.25%"
$ws.Range("D11").Value = "0.1248"
$ws.Range("E11").Value = "-8.04%"
$ws.Range("D12").Value = "0.1870"
$ws.Range("E12").Value = "-6.16%"
$ws.Range("D13").Value = "0.08801"
$ws.Range("E13").Value = "-3.60%"
$ws.Range("D14").Value = "0.03385"
$ws.Range("E14").Value = "-4.22%"
$ws.Range("D15").Value = "0.09691"
$ws.Range("E15").Value = "-1.30%"
$ws.Range("D16").Value = "0.001368"
$ws.Range("E16").Value = "-3.40%"
$ws.Range("D17").Value = "0.005964"
$ws.Range("E17").Value = "-3.49%"
$ws.Range("D18").Value = "3.588"
$ws.Range("E18").Value = "-2.72%"
$ws.Range("E19").Value = "-2.45%"
$ws.Range("D20").Value = "0.1270"
$ws.Range("E20").Value = "-4.10%"
$ws.Range("D21").Value = "5.015"
$ws.Range("E21").Value = "1.01%"
$ws.Range("D22").Value = "0.2489"
$ws.Range("E22").Value = "1.61%"
$ws.Range("D23").Value = "0.02104"
$ws.Range("E23").Value = "5,164.06%"
$ws.Range("D24").Value = "0.04326"
$ws.Range("E24").Value = "-0.69%"
$ws.Range("D25").Value = "0.001209"
$ws.Range("E25").Value = "-2.10%"
$ws.Range("D26").Value = "0.004232"
$ws.Range("E26").Value = "-12.00%"
$ws.Range("D27").Value = "0.0001349"
$ws.Range("E27").Value = "3.65%"
$ws.Range("D39").Value = "0.02174"
$ws.Range("E39").Value = "-6.18%"
$ws.Range("D40").Value = "0.04892"
$ws.Range("E40").Value = "-6.00%"
$ws.Range("D41").Value = "0.007701"
$ws.Range("E41").Value = "-0.58%"
$ws.Range("D42").Value = "0.009887"
$ws.Range("E42").Value = "-4.85%"
$ws.Range("D43").Value = "0.1339"
$ws.Range("E43").Value = "-5.18%"
$ws.Range("D44").Value = "0.001993"
$ws.Range("E44").Value = "-2.58%"
$ws.Range("D45").Value = "0.009840"
$ws.Range("E45").Value = "5.74%"
$ws.Range("D46").Value = "0.00006528"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.02%"
$ws.Range("D48").Value = "0.003000"
$ws.Range("E48").Value = "3.83%"
$ws.Range("D49").Value = "0.001300"
$ws.Range("E49").Value = "-23.05%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "0.02%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "0.02%"

foreach ($ref in $cellsToFormat) { $ws.Range($ref).Style = "Normal" }
